$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price (D) and volume-change (E) values, preserving original text formatting.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.533.85'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.78%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.645.12'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.71%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.25%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.005'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.28%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '303.23'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.04%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3825'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.15%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3598'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.83%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '50.94'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.96%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08166'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.08%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.223'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.20%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.006'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.21%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.30'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.03%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.422'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.09%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.399'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.53%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001216'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.12%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.652.34'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.75%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '97.32'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +2.48%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07031'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.15%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.756'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +2.93%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.45'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.74%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.005'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.21%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.65'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.69%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '23.530.80'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.73%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.477'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -2.85%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.016'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.02%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.22'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.86%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '153.35'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.99%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.232'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.60%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.73'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.47%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.836.84'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +1.53%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.037'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +8.54%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.263'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +5.58%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '12.04'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +4.16%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.053'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -2.55%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02793'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.62%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2494'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.03%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.08776'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.23%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.034'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.64%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06968'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.62%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.94'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +6.18%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6959'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.28%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.333'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.93%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.95'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +3.07%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6466'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.49%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.004'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.17%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.289'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.87%  '

$ws.Range("E48").Value = '  +0.34%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07849'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.05%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '127.64'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.87%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.171'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.15%  '
